# Apply cryptos list update (Mon Mar 11 17:21:33 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "72.409.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.35%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.036.63"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.22%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.736"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +20.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.760"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.37%  "
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000324"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.99"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +11.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.683.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.023.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.60%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.11"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +6.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").Value = "  -1.05%  "
$ws.Range("E19").Value = "  -1.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.183.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "441.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "104.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +17.91%  "
$ws.Range("E23").Value = "  +6.52%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "14.59"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.47"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.05"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.97%  "
$ws.Range("E29").Value = "  +2.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.56%  "
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "671.89"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +13.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "67.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "42.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.431"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0856"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.52"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +12.49%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.150"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0495"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.160"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +13.56%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.19%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.71"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.26%  "
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.02"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.26%  "
$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.00%  "
